$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2424.7144
$ws.Range("I15").Value = 2424.7144
$ws.Range("K15").Value = 7274.1432
$ws.Range("M15").Value = -7105.1432

$ws.Range("H40").Value = 890
$ws.Range("J40").Value = 890
$ws.Range("L40").Value = 890
$ws.Range("N40").Value = -1240

$ws.Range("H43").Value = 9776.5
$ws.Range("I43").Value = 9997.5
$ws.Range("J43").Value = 9555.5
$ws.Range("K43").Value = 9997.5
$ws.Range("L43").Value = 9555.5
$ws.Range("M43").Value = -9928.5
$ws.Range("N43").Value = -9693.5

$ws.Range("H107").Value = 43432.46
$ws.Range("I107").Value = 47238.453
$ws.Range("J107").Value = 22499.5
$ws.Range("K107").Value = 47238.453
$ws.Range("L107").Value = 22499.5
$ws.Range("M107").Value = -45318.453
$ws.Range("N107").Value = -26339.5

$ws.Range("H132").Value = 5299.2354
$ws.Range("I132").Value = 2175.0833
$ws.Range("J132").Value = 12797.2
$ws.Range("K132").Value = 6525.249899999999
$ws.Range("L132").Value = 38391.60000000001
$ws.Range("M132").Value = -3995.249899999999
$ws.Range("N132").Value = -43451.60000000001

$ws.Range("H137").Value = 11049.911
$ws.Range("I137").Value = 2393.6956
$ws.Range("K137").Value = 7181.0868
$ws.Range("M137").Value = -4631.0868

$ws.Range("H138").Value = 3176.8572
$ws.Range("J138").Value = 3742.5862
$ws.Range("L138").Value = 11227.7586
$ws.Range("N138").Value = -21507.7586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1340.5454
$ws.Range("I2").Value = 1522.36
$ws.Range("J2").Value = 772.375
$ws.Range("K2").Value = 1522.36
$ws.Range("L2").Value = 772.375
$ws.Range("M2").Value = -1409.36
$ws.Range("N2").Value = -998.375

$ws.Range("H32").Value = 4950.019
$ws.Range("I32").Value = 2257.2354
$ws.Range("J32").Value = 9768.684999999999
$ws.Range("K32").Value = 2257.2354
$ws.Range("L32").Value = 9768.684999999999
$ws.Range("M32").Value = -1970.2354
$ws.Range("N32").Value = -10342.685

$ws.Range("H110").Value = 5426.909
$ws.Range("I110").Value = 5426.909
$ws.Range("K110").Value = 5426.909
$ws.Range("M110").Value = -3381.909

$ws.Range("H116").Value = 1340.5454
$ws.Range("I116").Value = 1522.36
$ws.Range("J116").Value = 772.375
$ws.Range("K116").Value = 1522.36
$ws.Range("L116").Value = 772.375
$ws.Range("M116").Value = 771.6400000000001
$ws.Range("N116").Value = -5360.375

$ws.Range("H128").Value = 139999
$ws.Range("J128").Value = 139999
$ws.Range("L128").Value = 139999
$ws.Range("N128").Value = -149959

$ws.Range("H132").Value = 1906708.4
$ws.Range("I132").Value = 3283.1428
$ws.Range("J132").Value = 10789359
$ws.Range("K132").Value = 9849.428400000001
$ws.Range("L132").Value = 32368077
$ws.Range("M132").Value = -7319.428400000001
$ws.Range("N132").Value = -32373137

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1340.5454
$ws.Range("I3").Value = 1522.36
$ws.Range("J3").Value = 772.375
$ws.Range("K3").Value = 1522.36
$ws.Range("L3").Value = 772.375
$ws.Range("M3").Value = -1408.36
$ws.Range("N3").Value = -1000.375

$ws.Range("H20").Value = 20599.018
$ws.Range("I20").Value = 10193.6455
$ws.Range("J20").Value = 34039.293
$ws.Range("K20").Value = 10193.6455
$ws.Range("L20").Value = 34039.293
$ws.Range("M20").Value = -9946.645500000001
$ws.Range("N20").Value = -34533.293

$ws.Range("H80").Value = 1828.0625
$ws.Range("J80").Value = 2095
$ws.Range("L80").Value = 2095
$ws.Range("N80").Value = -4091

$ws.Range("H83").Value = 1828.0625
$ws.Range("J83").Value = 2095
$ws.Range("L83").Value = 10475
$ws.Range("N83").Value = -20459

$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws.Range("H133").Value = 149500
$ws.Range("I133").Value = 150000
$ws.Range("J133").Value = 149000
$ws.Range("K133").Value = 150000
$ws.Range("L133").Value = 149000
$ws.Range("M133").Value = -144940
$ws.Range("N133").Value = -159120

$ws.Range("H134").Value = 12462.78
$ws.Range("I134").Value = 7998.5405
$ws.Range("J134").Value = 53757
$ws.Range("K134").Value = 23995.6215
$ws.Range("L134").Value = 161271
$ws.Range("M134").Value = -21460.6215
$ws.Range("N134").Value = -166341

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 63743.207
$ws.Range("I31").Value = 127370
$ws.Range("J31").Value = 19481.088
$ws.Range("K31").Value = 127370
$ws.Range("L31").Value = 19481.088
$ws.Range("M31").Value = -127075
$ws.Range("N31").Value = -20071.088

$ws.Range("H34").Value = 63743.207
$ws.Range("I34").Value = 127370
$ws.Range("J34").Value = 19481.088
$ws.Range("K34").Value = 127370
$ws.Range("L34").Value = 19481.088
$ws.Range("M34").Value = -127168
$ws.Range("N34").Value = -19885.088

$ws.Range("H109").Value = 17666.666
$ws.Range("J109").Value = 17666.666
$ws.Range("L109").Value = 17666.666
$ws.Range("N109").Value = -19746.666

$ws.Range("H132").Value = 1479374.1
$ws.Range("I132").Value = 3859.5334
$ws.Range("K132").Value = 11578.6002
$ws.Range("M132").Value = -9048.600199999999

$ws.Range("H134").Value = 7661.394
$ws.Range("I134").Value = 3437.3333
$ws.Range("K134").Value = 10311.9999
$ws.Range("M134").Value = -7776.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 493.4
$ws.Range("I5").Value = 149
$ws.Range("J5").Value = 531.6667
$ws.Range("K5").Value = 447
$ws.Range("L5").Value = 1595.0001
$ws.Range("M5").Value = -335
$ws.Range("N5").Value = -1819.0001

$ws.Range("H80").Value = 35999.668
$ws.Range("I80").Value = 29999.5
$ws.Range("K80").Value = 89998.5
$ws.Range("M80").Value = -89062.5

$ws.Range("H83").Value = 35999.668
$ws.Range("I83").Value = 29999.5
$ws.Range("K83").Value = 269995.5
$ws.Range("M83").Value = -265315.5

$ws.Range("H114").Value = 1221.2
$ws.Range("I114").Value = 628
$ws.Range("K114").Value = 1884
$ws.Range("M114").Value = 1370

$ws.Range("H135").Value = 493.4
$ws.Range("I135").Value = 149
$ws.Range("J135").Value = 531.6667
$ws.Range("K135").Value = 1341
$ws.Range("L135").Value = 4785.0003
$ws.Range("M135").Value = 1194
$ws.Range("N135").Value = -9855.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12228.161
$ws.Range("I70").Value = 11322.723
$ws.Range("K70").Value = 11322.723
$ws.Range("M70").Value = -11052.723

$ws.Range("H73").Value = 12228.161
$ws.Range("I73").Value = 11322.723
$ws.Range("K73").Value = 11322.723
$ws.Range("M73").Value = -10386.723

$ws.Range("H132").Value = 715009.8
$ws.Range("I132").Value = 5338
$ws.Range("K132").Value = 16014
$ws.Range("M132").Value = -13484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2760.875
$ws.Range("J16").Value = 2664.3333
$ws.Range("L16").Value = 2664.3333
$ws.Range("N16").Value = -3004.3333

$ws.Range("H22").Value = 1400.3334
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 1750.5
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 1750.5
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -2340.5

$ws.Range("H27").Value = 1400.3334
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 1750.5
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 1750.5
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -1964.5

$ws.Range("H46").Value = 3386.3157
$ws.Range("I46").Value = 787
$ws.Range("K46").Value = 787
$ws.Range("M46").Value = -599

$ws.Range("H109").Value = 12500
$ws.Range("J109").Value = 12500
$ws.Range("L109").Value = 12500
$ws.Range("N109").Value = -15274

$ws.Range("H122").Value = 9023.388999999999
$ws.Range("I122").Value = 12401
$ws.Range("K122").Value = 37203
$ws.Range("M122").Value = -34753

$ws.Range("H130").Value = 51999.75
$ws.Range("J130").Value = 53999.5
$ws.Range("L130").Value = 53999.5
$ws.Range("N130").Value = -64039.5

$ws.Range("H136").Value = 1036135.6
$ws.Range("I136").Value = 18553.77
$ws.Range("J136").Value = 1771055.9
$ws.Range("K136").Value = 55661.31
$ws.Range("L136").Value = 5313167.699999999
$ws.Range("M136").Value = -53111.31
$ws.Range("N136").Value = -5318267.699999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8920.5
$ws.Range("I107").Value = 9640.272000000001
$ws.Range("J107").Value = 1003
$ws.Range("K107").Value = 28920.816
$ws.Range("L107").Value = 3009
$ws.Range("M107").Value = -27000.816
$ws.Range("N107").Value = -6849

$ws.Range("H109").Value = 22666.666
$ws.Range("J109").Value = 22666.666
$ws.Range("L109").Value = 22666.666
$ws.Range("N109").Value = -25440.666

$ws.Range("H122").Value = 3870.6206
$ws.Range("I122").Value = 2587.75
$ws.Range("J122").Value = 6721.4443
$ws.Range("K122").Value = 7763.25
$ws.Range("L122").Value = 20164.3329
$ws.Range("M122").Value = -5313.25
$ws.Range("N122").Value = -25064.3329

$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -54820

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
